{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the body currently ends with \" profile card\"\n// followed by the _GoBack bookmark. We append a brand-new paragraph\n// after it that holds the fontawesome.io hyperlink plus a trailing space,\n// matching the target diff.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst newParagraph = lastParagraph.insertParagraph(\"\", \"After\");\n\n// Insert the hyperlink run with its own text, then mark that sub-range as\n// a hyperlink (gives it the \"Hyperlink\" style + the external address).\nconst linkRange = newParagraph.insertText(\"http://fontawesome.io/examples/#\", \"Start\");\nlinkRange.hyperlink = \"http://fontawesome.io/examples/#\";\n\n// Trailing run: a single space, appended right after the hyperlink text.\nconst spaceRange = newParagraph.insertText(\" \", \"End\");\n\n// The \"_GoBack\" bookmark used to sit at the end of the old last paragraph\n// (right after \" profile card\"); it now belongs at the end of the new\n// trailing paragraph, after the appended space.\ncontext.document.deleteBookmark(\"_GoBack\");\nspaceRange.getRange(\"End\").insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Append a brand-new paragraph after the document's current last paragraph\n# (the one ending \" profile card\"), holding the fontawesome.io hyperlink\n# followed by a trailing space - mirrors the target diff.\n$lastParagraph = $d.Paragraphs.Last\n$end = $lastParagraph.Range\n$end.Collapse(0)  # wdCollapseEnd\n$end.InsertParagraphAfter()\n\n$newRange = $d.Paragraphs.Last.Range\n$newRange.Collapse(0)\n\n$d.Hyperlinks.Add($newRange, \"http://fontawesome.io/examples/#\") | Out-Null\n\n$spaceRange = $d.Paragraphs.Last.Range\n$spaceRange.Collapse(0)\n$spaceRange.InsertAfter(\" \")\n\n# The \"_GoBack\" bookmark used to mark the end of the old last paragraph;\n# now it belongs at the end of the new trailing paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n$bookmarkRange = $d.Paragraphs.Last.Range\n$bookmarkRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange) | Out-Null\n"}
